$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.235.91"
$ws.Range("E2").Value = "'  -1.57%  "

$ws.Range("D3").Value = "'1.850.52"
$ws.Range("E3").Value = "'  -2.62%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.05%  "

$ws.Range("D5").Value = "'233.25"
$ws.Range("E5").Value = "'  -2.44%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  +0.03%  "

$ws.Range("D7").Value = "'0.4698"
$ws.Range("E7").Value = "'  -2.40%  "

$ws.Range("D8").Value = "'0.2705"
$ws.Range("E8").Value = "'  -4.89%  "

$ws.Range("D9").Value = "'0.06375"
$ws.Range("E9").Value = "'  -2.76%  "

$ws.Range("D10").Value = "'1.860.69"
$ws.Range("E10").Value = "'  -2.88%  "

$ws.Range("D11").Value = "'0.07424"
$ws.Range("E11").Value = "'  -0.54%  "

$ws.Range("D12").Value = "'16.26"
$ws.Range("E12").Value = "'  -2.87%  "

$ws.Range("D13").Value = "'4.941"
$ws.Range("E13").Value = "'  -3.53%  "

$ws.Range("D14").Value = "'85.08"
$ws.Range("E14").Value = "'  -3.56%  "

$ws.Range("D15").Value = "'0.6265"
$ws.Range("E15").Value = "'  -6.09%  "

$ws.Range("D16").Value = "'30.170.97"
$ws.Range("E16").Value = "'  -1.69%  "

$ws.Range("D17").Value = "'0.9999"
$ws.Range("E17").Value = "'  +0.02%  "

$ws.Range("D18").Value = "'229.14"
$ws.Range("E18").Value = "'  -1.13%  "

$ws.Range("D19").Value = "'12.58"
$ws.Range("E19").Value = "'  -5.64%  "

$ws.Range("D20").Value = "'0.000007342"
$ws.Range("E20").Value = "'  -3.71%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.090.76"
$ws.Range("E21").Value = "'  -4.43%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "'  -0.03%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.935"
$ws.Range("E23").Value = "'  -7.01%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.927"
$ws.Range("E24").Value = "'  -5.03%  "

$ws.Range("D25").Value = "'166.78"
$ws.Range("E25").Value = "'  -2.01%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.221"
$ws.Range("E26").Value = "'  -1.07%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'17.81"
$ws.Range("E27").Value = "'  -4.85%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.871"
$ws.Range("E28").Value = "'  -4.98%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.1034"
$ws.Range("E29").Value = "'  +2.87%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.381"
$ws.Range("E30").Value = "'  -1.46%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.119"
$ws.Range("E31").Value = "'  -5.64%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.875"
$ws.Range("E32").Value = "'  -4.23%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04887"
$ws.Range("E33").Value = "'  -4.03%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.152"
$ws.Range("E34").Value = "'  -5.46%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7091"
$ws.Range("E35").Value = "'  -6.45%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'0.9998"
$ws.Range("E36").Value = "'  +0.01%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.695"
$ws.Range("E37").Value = "'  -0.71%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01862"
$ws.Range("E38").Value = "'  -1.03%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.632"
$ws.Range("E39").Value = "'  -0.94%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9035"
$ws.Range("E40").Value = "'  -1.82%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.944"
$ws.Range("E41").Value = "'  -6.90%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'105.13"
$ws.Range("E42").Value = "'  -1.76%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9981"
$ws.Range("E43").Value = "'  -0.65%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.552"
$ws.Range("E44").Value = "'  -3.78%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4083"
$ws.Range("E45").Value = "'  -5.22%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.055"
$ws.Range("E46").Value = "'  -5.19%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'60.24"
$ws.Range("E47").Value = "'  -6.40%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1190"
$ws.Range("E48").Value = "'  -6.74%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.643"
$ws.Range("E49").Value = "'  -3.55%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'33.10"
$ws.Range("E50").Value = "'  -2.45%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.387"
$ws.Range("E51").Value = "'  -6.86%  "
